# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn and de-de sheets get their "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns populated, with a new hyperlink on the
#    "Latest Target File" cell.
#  - A few columns are widened to fit the new, longer content.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b223eb013c257fb9e2351bf965bcf4f75fa51b97/e2e/a.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b223eb013c257fb9e2351bf965bcf4f75fa51b97/e2e/b.md"

# Column width that round-trips (through this engine's 1/6-character
# quantization) to the widened width used by the real workbook (~29.98 -> 30).
$wideColWidth = 29.166666666666668

### --- Overview sheet --- ###
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

### --- zh-cn sheet --- ###
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("I3").Value = "a.md"

$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-22 04:47:28"
$wsZh.Range("K3").Value = "2016-08-22 04:47:28"

$wsZh.Columns.Item(3).ColumnWidth = $wideColWidth
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# Rebuild hyperlinks so the new "Latest Target File" links land in the same
# order/ids as the original workbook: A2, I2, A3, I3
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlB, [Type]::Missing, [Type]::Missing, "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")

### --- de-de sheet --- ###
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("I3").Value = "a.md"

$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-22 04:47:35"
$wsDe.Range("K3").Value = "2016-08-22 04:47:35"

$wsDe.Columns.Item(3).ColumnWidth = $wideColWidth
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlB, [Type]::Missing, [Type]::Missing, "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
